# Add 2022-Q3 data
# --------------------------------------------------------------------------
# 1) Build the new "2022-Q3" worksheet by copying the existing "2022-Q2"
#    sheet (so it inherits identical layout/column widths/styles), placing
#    it right before "2022-Q2" (i.e. as the 2nd sheet overall).
# --------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Overwrite the data rows (2-6) of the new sheet with the 2022-Q3 figures.
# Columns B-G are stored as text (fund codes / percentages kept as text in
# the source data), column A/H stay numeric.
$q3Rows = @(
    @("590001", "中邮核心优选混合",         "12.91", "81.83", "4.92", "0.6352", 2),
    @("009686", "华夏磐利一年定期开放混合A", "10.76", "64.78", "2.48", "0.2668", 9),
    @("010965", "中银鑫新消费成长混合A",     "3.22",  "86.59", "2.32", "0.0747", 10),
    @("010962", "中银鑫新消费成长混合C",     "0.67",  "86.59", "2.32", "0.0155", 10),
    @("009687", "华夏磐利一年定期开放混合C", "0.43",  "64.78", "2.48", "0.0107", 9)
)

for ($i = 0; $i -lt $q3Rows.Count; $i++) {
    $r = $i + 2
    $row = $q3Rows[$i]
    $q3.Cells.Item($r, 1).Value = $i
    $q3.Cells.Item($r, 2).Value = "'" + $row[0]
    $q3.Cells.Item($r, 3).Value = $row[1]
    $q3.Cells.Item($r, 4).Value = "'" + $row[2]
    $q3.Cells.Item($r, 5).Value = "'" + $row[3]
    $q3.Cells.Item($r, 6).Value = "'" + $row[4]
    $q3.Cells.Item($r, 7).Value = "'" + $row[5]
    $q3.Cells.Item($r, 8).Value = $row[6]
}

# --------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert the 2022-Q3 row at the top
#    of the data (row 2), pushing the other quarters down by one row.
# --------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Range("A2:D5").Clear()

$summaryRows = @(
    @("2022-Q3", 5, 1),
    @("2022-Q2", 5, 1.78),
    @("2022-Q1", 2, 0.66),
    @("2021-Q4", 5, 2.14),
    @("2020-Q4", 2, 0.09)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summary.Cells.Item($r, 1).Value = $i
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
}

# Re-apply the bold/centered/bordered style used by column A (same style as
# the header cells) which got lost after Clear().
$summary.Range("B1").Copy()
$summary.Range("A2:A6").PasteSpecial(-4122)

# --------------------------------------------------------------------------
# 3) Restore "2020-Q4" as the active/selected sheet (it was the last sheet
#    and tab-selected before this edit).
# --------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
